$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.023.69"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "'1.823.48"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "'309.43"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  -1.29%  "

$ws.Range("D8").Value = "'0.3661"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("D9").Value = "'0.07240"
$ws.Range("E9").Value = "  -2.90%  "

$ws.Range("D10").Value = "'0.8609"
$ws.Range("E10").Value = "  -2.71%  "

$ws.Range("D11").Value = "'19.86"
$ws.Range("E11").Value = "  -2.80%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07539"
$ws.Range("E12").Value = "  +2.80%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.737.36"
$ws.Range("E13").Value = "  -7.35%  "

$ws.Range("D14").Value = "'5.331"
$ws.Range("E14").Value = "  -2.05%  "

$ws.Range("D15").Value = "'92.04"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "'6.510"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").Value = "'0.000008644"
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "'27.082.36"
$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'14.49"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").Value = "'5.149"
$ws.Range("E22").Value = "  -2.79%  "

$ws.Range("D23").Value = "'10.52"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").Value = "'1.891.85"
$ws.Range("E24").Value = "  -9.72%  "

$ws.Range("D25").Value = "'151.57"
$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").Value = "'1.842"
$ws.Range("E26").Value = "  -2.65%  "

$ws.Range("D27").Value = "'18.17"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").Value = "'2.070"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("D29").Value = "'5.123"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").Value = "'115.35"
$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("D31").Value = "'0.08881"
$ws.Range("E31").Value = "  -1.34%  "

$ws.Range("D32").Value = "'2.960"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").Value = "'4.425"
$ws.Range("E33").Value = "  -2.62%  "

$ws.Range("D34").Value = "'1.132"
$ws.Range("E34").Value = "  -3.80%  "

$ws.Range("D35").Value = "'0.7201"
$ws.Range("E35").Value = "  -4.49%  "

$ws.Range("D36").Value = "'1.080"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").Value = "'0.05260"
$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.415"
$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01924"
$ws.Range("E39").Value = "  -1.48%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.931"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("D41").Value = "'7.155"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").Value = "'0.5163"
$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1629"
$ws.Range("E43").Value = "  -1.78%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.8576"
$ws.Range("E44").Value = "  -15.12%  "

$ws.Range("D45").Value = "'8.185"
$ws.Range("E45").Value = "  -3.43%  "

$ws.Range("D46").Value = "'0.4822"
$ws.Range("E46").Value = "  -1.84%  "

$ws.Range("D47").Value = "'1.006"
$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("D48").Value = "'10.12"
$ws.Range("E48").Value = "  -4.07%  "

$ws.Range("D49").Value = "'103.06"
$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").Value = "'1.623"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("D51").Value = "'0.06249"
$ws.Range("E51").Value = "  -0.77%  "
